# Fruta / hortaliza, semanal
# Reorders the per-row "Fecha" (D), "Volumen" (M), "Precio minimo" (N),
# "Precio maximo" (O), "Precio promedio ponderado" (P) and "Precio $/Kg" (S)
# values across rows 2-9 to reflect the corrected weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44186
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

$ws.Range("D3").Value = 44188
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 3000

$ws.Range("D4").Value = 44193
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("S4").Value = 3000

$ws.Range("D5").Value = 44175
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 4000

$ws.Range("D6").Value = 44179
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 4000

$ws.Range("D7").Value = 44189
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range("S7").Value = 3000

$ws.Range("D8").Value = 44196
$ws.Range("M8").Value = 56
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 3000

$ws.Range("D9").Value = 44181
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000
